# Update "paises.xlsx" country data / provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Julio de 2020 a las 12:05"

# --- Row 6: India -- refreshed numbers ---
$ws.Range("B6").Value = 1241654
$ws.Range("C6").Value = 1970
$ws.Range("D6").Value = 784440
$ws.Range("E6").Value = 427308
$ws.Range("G6").Value = 16
$ws.Range("H6").Value = 29906

# --- Row 33: Filipinas -- refreshed numbers ---
$ws.Range("B33").Value = 74390
$ws.Range("C33").Value = 2200
$ws.Range("D33").Value = 24383
$ws.Range("E33").Value = 48136
$ws.Range("G33").Value = 28
$ws.Range("H33").Value = 1871

# --- Row 34: Oman -- refreshed numbers ---
$ws.Range("B34").Value = 72646
$ws.Range("C34").Value = 1099
$ws.Range("D34").Value = 51349
$ws.Range("E34").Value = 20942
$ws.Range("G34").Value = 6
$ws.Range("H34").Value = 355

# --- Rows 48-49: Rumania moves above Guatemala (ranking swap) ---
# Row 48 becomes Rumania with newly refreshed numbers
$ws.Range("A48").Value = "Rumania"
$ws.Range("B48").Value = 41275
$ws.Range("C48").Value = 1112
$ws.Range("D48").Value = 24862
$ws.Range("E48").Value = 14287
$ws.Range("G48").Value = 25
$ws.Range("H48").Value = 2126

# Row 49 becomes Guatemala, keeping its previous (unchanged) numbers
$ws.Range("A49").Value = "Guatemala"
$ws.Range("B49").Value = 41135
$ws.Range("C49").Value = 0
$ws.Range("D49").Value = 27756
$ws.Range("E49").Value = 11806
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 1573

# --- Rows 73-76: El Salvador jumps above Dinamarca / Australia / Venezuela ---
# Row 73 becomes El Salvador with newly refreshed numbers
$ws.Range("A73").Value = "El Salvador"
$ws.Range("B73").Value = 13377
$ws.Range("C73").Value = 402
$ws.Range("D73").Value = 7276
$ws.Range("E73").Value = 5729
$ws.Range("G73").Value = 9
$ws.Range("H73").Value = 372

# Row 74 becomes Dinamarca, keeping its previous (unchanged) numbers
$ws.Range("A74").Value = "Dinamarca"
$ws.Range("B74").Value = 13350
$ws.Range("C74").Value = 0
$ws.Range("D74").Value = 12274
$ws.Range("E74").Value = 465
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 611

# Row 75 becomes Australia, keeping its previous (unchanged) numbers
$ws.Range("A75").Value = "Australia"
$ws.Range("B75").Value = 13302
$ws.Range("C75").Value = 406
$ws.Range("D75").Value = 8656
$ws.Range("E75").Value = 4513
$ws.Range("G75").Value = 5
$ws.Range("H75").Value = 133

# Row 76 becomes Venezuela, keeping its previous (unchanged) numbers
$ws.Range("A76").Value = "Venezuela"
$ws.Range("B76").Value = 13164
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 7471
$ws.Range("E76").Value = 5569
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 124

# --- Row 85: Noruega -- refreshed numbers ---
$ws.Range("B85").Value = 9062
$ws.Range("C85").Value = 3
$ws.Range("D85").Value = 8674
$ws.Range("E85").Value = 133

# --- Row 86: Malasia -- refreshed numbers ---
$ws.Range("B86").Value = 8840
$ws.Range("C86").Value = 9
$ws.Range("D86").Value = 8574
$ws.Range("E86").Value = 143

# --- Row 89: Finlandia -- refreshed numbers ---
$ws.Range("B89").Value = 7372
$ws.Range("C89").Value = 10
$ws.Range("E89").Value = 124

# --- Row 124: Eslovenia -- refreshed numbers ---
$ws.Range("D124").Value = 1661
$ws.Range("E124").Value = 257
